$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 2 "Bitcoin"
Set-TextValue 2 3 "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextValue 2 4 "29.341.81"
Set-TextValue 2 5 "  -0.24%  "

Set-TextValue 3 2 "Ethereum"
Set-TextValue 3 3 "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextValue 3 4 "1.846.60"
Set-TextValue 3 5 "  -0.18%  "

Set-TextValue 4 2 "TetherUSD"
Set-TextValue 4 3 "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextValue 4 4 "0.9980"
Set-TextValue 4 5 "  -0.16%  "

Set-TextValue 5 2 "BNB"
Set-TextValue 5 3 "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue 5 4 "240.20"
Set-TextValue 5 5 "  -0.26%  "

Set-TextValue 6 2 "XRP"
Set-TextValue 6 3 "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue 6 4 "0.6267"
Set-TextValue 6 5 "  -0.55%  "

Set-TextValue 7 2 "USDC"
Set-TextValue 7 3 "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue 7 4 "0.9991"
Set-TextValue 7 5 "  -0.12%  "

Set-TextValue 8 2 "Dogecoin"
Set-TextValue 8 3 "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue 8 4 "0.07602"
Set-TextValue 8 5 "  -1.00%  "

Set-TextValue 9 2 "Cardano"
Set-TextValue 9 3 "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue 9 4 "0.2903"
Set-TextValue 9 5 "  -1.27%  "

Set-TextValue 10 2 "Solana"
Set-TextValue 10 3 "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue 10 4 "24.65"
Set-TextValue 10 5 "  +0.53%  "

Set-TextValue 11 2 "TRON"
Set-TextValue 11 3 "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue 11 4 "0.07737"
Set-TextValue 11 5 "  -0.11%  "

Set-TextValue 12 2 "Polkadot"
Set-TextValue 12 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue 12 4 "5.021"
Set-TextValue 12 5 "  -0.04%  "

Set-TextValue 13 2 "Polygon"
Set-TextValue 13 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue 13 4 "0.6781"
Set-TextValue 13 5 "  -0.46%  "

Set-TextValue 14 2 "ShibaInu"
Set-TextValue 14 3 "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue 14 4 "0.00001061"
Set-TextValue 14 5 "  -3.04%  "

Set-TextValue 15 2 "Litecoin"
Set-TextValue 15 3 "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue 15 4 "82.91"
Set-TextValue 15 5 "  -0.94%  "

Set-TextValue 16 2 "Uniswap"
Set-TextValue 16 3 "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue 16 4 "6.126"
Set-TextValue 16 5 "  -0.47%  "

Set-TextValue 17 2 "WrappedBTC"
Set-TextValue 17 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue 17 4 "29.373.84"
Set-TextValue 17 5 "  -0.27%  "

Set-TextValue 18 2 "BitcoinCash"
Set-TextValue 18 3 "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue 18 4 "228.17"
Set-TextValue 18 5 "  -0.48%  "

Set-TextValue 19 2 "Avalanche"
Set-TextValue 19 3 "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue 19 4 "12.34"
Set-TextValue 19 5 "  -1.17%  "

Set-TextValue 20 2 "Dai"
Set-TextValue 20 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue 20 4 "0.9989"
Set-TextValue 20 5 "  -0.12%  "

Set-TextValue 21 2 "Chainlink"
Set-TextValue 21 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue 21 4 "7.487"
Set-TextValue 21 5 "  +0.45%  "

Set-TextValue 22 2 "BinanceUSD"
Set-TextValue 22 3 "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue 22 4 "0.9992"
Set-TextValue 22 5 "  -0.12%  "

Set-TextValue 23 2 "Monero"
Set-TextValue 23 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 23 4 "158.54"
Set-TextValue 23 5 "  +0.81%  "

Set-TextValue 24 2 "Stellar"
Set-TextValue 24 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue 24 4 "0.1382"
Set-TextValue 24 5 "  -0.38%  "

Set-TextValue 25 2 "Cosmos"
Set-TextValue 25 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue 25 4 "8.434"
Set-TextValue 25 5 "  +0.41%  "

Set-TextValue 26 2 "EthereumClassic"
Set-TextValue 26 3 "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue 26 4 "17.66"
Set-TextValue 26 5 "  -0.10%  "

Set-TextValue 27 2 "Toncoin"
Set-TextValue 27 3 "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue 27 4 "1.438"
Set-TextValue 27 5 "  +9.77%  "

Set-TextValue 28 2 "PancakeSwap"
Set-TextValue 28 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue 28 4 "1.467"
Set-TextValue 28 5 "  +0.06%  "

Set-TextValue 29 2 "Hedera"
Set-TextValue 29 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 29 4 "0.05606"
Set-TextValue 29 5 "  -1.87%  "

Set-TextValue 30 2 "Filecoin"
Set-TextValue 30 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue 30 4 "4.098"
Set-TextValue 30 5 "  -0.51%  "

Set-TextValue 31 2 "InternetComputer(DFINITY)"
Set-TextValue 31 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue 31 4 "4.064"
Set-TextValue 31 5 "  +0.25%  "

Set-TextValue 32 2 "ARBITRUM"
Set-TextValue 32 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue 32 4 "1.160"
Set-TextValue 32 5 "  -0.19%  "

Set-TextValue 33 2 "LidoDAOToken"
Set-TextValue 33 3 "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue 33 4 "1.830"
Set-TextValue 33 5 "  -1.08%  "

Set-TextValue 34 2 "ImmutableX"
Set-TextValue 34 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 34 4 "0.6973"
Set-TextValue 34 5 "  -1.51%  "

Set-TextValue 35 2 "HuobiToken"
Set-TextValue 35 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue 35 4 "2.583"
Set-TextValue 35 5 "  -0.17%  "

Set-TextValue 36 2 "Maker"
Set-TextValue 36 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue 36 4 "1.232.71"
Set-TextValue 36 5 "  +0.79%  "

Set-TextValue 37 2 "VeChain"
Set-TextValue 37 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 37 4 "0.01797"
Set-TextValue 37 5 "  -0.12%  "

Set-TextValue 38 2 "MXToken"
Set-TextValue 38 3 "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue 38 4 "2.727"
Set-TextValue 38 5 "  -1.68%  "

Set-TextValue 39 2 "FraxShare"
Set-TextValue 39 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue 39 4 "6.354"
Set-TextValue 39 5 "  -1.83%  "

Set-TextValue 40 2 "TrustWalletToken"
Set-TextValue 40 3 "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue 40 4 "0.8985"
Set-TextValue 40 5 "  -1.14%  "

Set-TextValue 41 2 "PaxDollar"
Set-TextValue 41 3 "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue 41 4 "0.9990"
Set-TextValue 41 5 "  -0.16%  "

Set-TextValue 42 2 "Quant"
Set-TextValue 42 3 "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue 42 4 "101.44"
Set-TextValue 42 5 "  -0.12%  "

Set-TextValue 43 2 "Aave"
Set-TextValue 43 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue 43 4 "65.32"
Set-TextValue 43 5 "  -1.37%  "

Set-TextValue 44 2 "Aptos"
Set-TextValue 44 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue 44 4 "7.190"
Set-TextValue 44 5 "  +0.63%  "

Set-TextValue 45 2 "BabyDogeCoin"
Set-TextValue 45 3 "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue 45 4 "0.00000000117"
Set-TextValue 45 5 "  -3.06%  "

Set-TextValue 46 2 "TheSandbox"
Set-TextValue 46 3 "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue 46 4 "0.3990"
Set-TextValue 46 5 "  -0.64%  "

Set-TextValue 47 2 "EnergySwap"
Set-TextValue 47 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue 47 4 "8.988"
Set-TextValue 47 5 "  -0.07%  "

Set-TextValue 48 2 "RenderToken"
Set-TextValue 48 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 48 4 "1.684"
Set-TextValue 48 5 "  -0.20%  "

Set-TextValue 49 2 "Algorand"
Set-TextValue 49 3 "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue 49 4 "0.1141"
Set-TextValue 49 5 "  +1.15%  "

Set-TextValue 50 2 "Cronos"
Set-TextValue 50 3 "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue 50 4 "0.05695"
Set-TextValue 50 5 "  -0.28%  "

Set-TextValue 51 2 "Mantle"
Set-TextValue 51 3 "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue 51 4 "0.4622"
Set-TextValue 51 5 "  -0.09%  "

Write-Output "done"